$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (03-dec) before column EF (01-oct.) ---
$ws1 = $wb.Worksheets.Item(1)

# Column EF is the 136th column; inserting there shifts 01-oct..31-oct (EF..FJ) to EG..FK
$ws1.Columns.Item(136).Insert()

# Header for the newly inserted column
$ws1.Range("EF1").Value = "03-dec"

# Fill the new column's data rows (2-25) with the placeholder used for missing data
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 136).Value = "-"
}

# --- Sheet "Gaz": append a new row with the latest date/price ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A166").NumberFormat = "@"
$ws2.Range("A166").Value = "2025-12-01"
$ws2.Range("A166").NumberFormat = "General"
$ws2.Range("A166").Style = "Normal"
$ws2.Range("B166").Value = 27.2

# --- Sheet "CO2": append a new row with the latest date/price ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A166").NumberFormat = "@"
$ws3.Range("A166").Value = "2025-12-01"
$ws3.Range("A166").NumberFormat = "General"
$ws3.Range("A166").Style = "Normal"
$ws3.Range("B166").Value = 82.68000000000001
